$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some of which, e.g. "277.46",
# "0.3615", look like plain numbers to Excel's auto-detection).
# Force a Text format before assigning so the new value is kept
# as a literal string (preserving exact digits/trailing zeros),
# then clear the format again so the cell style matches the
# original (no explicit style on these data cells).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '20.540.40'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.472.25'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  +4.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '277.18'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3613'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3074'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.64'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.071'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06650'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.515'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.14'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.170'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9582'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001026'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.473.49'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05932'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.75'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.488'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.52'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.15'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.267'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '20.545.35'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.91'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.127'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.64%  '
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.633.34'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.86'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.900'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.955'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.08014'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.8021'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  +4.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.218'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.27%  '
$ws.Range('E37').Value = '  -2.79%  '
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02055'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9587'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.38'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1870'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.403'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5270'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.21'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.42'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5202'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.814'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06462'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9868'
$ws.Range('D51').ClearFormats()
